$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.193.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4666"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.097"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6727"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.185.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.529"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.112.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007283"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.197"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.918"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.376"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09715"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.418"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.473"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.094"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.66%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  +1.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7048"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.78%  "

$ws.Range("E37").Value = "  +0.55%  "

$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.239"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.942"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8447"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("E46").Value = "  -0.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.175"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.173"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "934.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("E51").Value = "  -3.00%  "
